$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings) - new crime-report week.
# ---------------------------------------------------------------------------
$ws.Cells.Item(8,1).Value = "Volume 32   Number  22"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# ---------------------------------------------------------------------------
# Helper number formats matching the workbook's existing styles:
#   style 14 -> "#,##0"                 (plain integer count)
#   style 15 -> '#,##0.0;"-"#,##0.0'    (one-decimal percent-change)
# ---------------------------------------------------------------------------
$fmtCount = "#,##0"
$fmtPct   = '#,##0.0;"-"#,##0.0'

# Cells that are currently text placeholders ("0" / "***.*") and must become
# real numbers - give them the right NumberFormat before writing the value so
# the engine maps them back onto the workbook's existing numeric styles.
$ws.Cells.Item(15,3).NumberFormat  = $fmtCount
$ws.Cells.Item(15,6).NumberFormat  = $fmtCount
$ws.Cells.Item(17,3).NumberFormat  = $fmtCount
$ws.Cells.Item(18,3).NumberFormat  = $fmtCount
$ws.Cells.Item(20,4).NumberFormat  = $fmtCount
$ws.Cells.Item(20,5).NumberFormat  = $fmtPct
$ws.Cells.Item(23,4).NumberFormat  = $fmtCount
$ws.Cells.Item(23,5).NumberFormat  = $fmtPct
$ws.Cells.Item(27,3).NumberFormat  = $fmtCount
$ws.Cells.Item(27,6).NumberFormat  = $fmtCount
$ws.Cells.Item(29,4).NumberFormat  = $fmtCount
$ws.Cells.Item(29,5).NumberFormat  = $fmtPct
$ws.Cells.Item(30,4).NumberFormat  = $fmtCount
$ws.Cells.Item(30,5).NumberFormat  = $fmtPct

# C20 is the only cell going the other way (number -> text placeholder "0").
# Force text storage via a Text number format, then restore the canonical
# "General" text style (13) by copying formats from a cell that already
# carries it (the row label is always style 13 and never changes).
$ws.Cells.Item(20,3).NumberFormat = "@"
$ws.Cells.Item(20,3).Value = "0"
$ws.Cells.Item(20,1).Copy() | Out-Null
$ws.Cells.Item(20,3).PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Cells.Item(15,3).Value  = 1
$ws.Cells.Item(15,6).Value  = 1
$ws.Cells.Item(15,9).Value  = 2
$ws.Cells.Item(15,11).Value = 0
$ws.Cells.Item(15,12).Value = 100
$ws.Cells.Item(15,13).Value = -33.333333333333
$ws.Cells.Item(15,14).Value = -60

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Cells.Item(16,4).Value  = 1
$ws.Cells.Item(16,5).Value  = 0
$ws.Cells.Item(16,6).Value  = 4
$ws.Cells.Item(16,8).Value  = -42.857142857142
$ws.Cells.Item(16,9).Value  = 18
$ws.Cells.Item(16,10).Value = 24
$ws.Cells.Item(16,11).Value = -25
$ws.Cells.Item(16,12).Value = 20
$ws.Cells.Item(16,13).Value = -35.714285714285
$ws.Cells.Item(16,14).Value = -83.636363636363

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Cells.Item(17,3).Value  = 1
$ws.Cells.Item(17,4).Value  = 7
$ws.Cells.Item(17,5).Value  = -85.714285714285
$ws.Cells.Item(17,6).Value  = 9
$ws.Cells.Item(17,7).Value  = 21
$ws.Cells.Item(17,8).Value  = -57.142857142857
$ws.Cells.Item(17,9).Value  = 36
$ws.Cells.Item(17,10).Value = 71
$ws.Cells.Item(17,11).Value = -49.295774647887
$ws.Cells.Item(17,12).Value = 5.882352941176
$ws.Cells.Item(17,13).Value = 2.857142857142
$ws.Cells.Item(17,14).Value = -55.555555555555

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Cells.Item(18,3).Value  = 2
$ws.Cells.Item(18,5).Value  = 100
$ws.Cells.Item(18,6).Value  = 5
$ws.Cells.Item(18,7).Value  = 4
$ws.Cells.Item(18,8).Value  = 25
$ws.Cells.Item(18,9).Value  = 22
$ws.Cells.Item(18,10).Value = 22
$ws.Cells.Item(18,11).Value = 0
$ws.Cells.Item(18,12).Value = 57.142857142857
$ws.Cells.Item(18,14).Value = -89.215686274509

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Cells.Item(19,3).Value  = 6
$ws.Cells.Item(19,4).Value  = 8
$ws.Cells.Item(19,5).Value  = -25
$ws.Cells.Item(19,6).Value  = 12
$ws.Cells.Item(19,7).Value  = 14
$ws.Cells.Item(19,8).Value  = -14.285714285714
$ws.Cells.Item(19,9).Value  = 50
$ws.Cells.Item(19,10).Value = 62
$ws.Cells.Item(19,11).Value = -19.354838709677
$ws.Cells.Item(19,12).Value = -26.470588235294
$ws.Cells.Item(19,13).Value = 13.636363636363
$ws.Cells.Item(19,14).Value = -20.634920634920

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Cells.Item(20,4).Value  = 1
$ws.Cells.Item(20,5).Value  = -100
$ws.Cells.Item(20,6).Value  = 2
$ws.Cells.Item(20,7).Value  = 2
$ws.Cells.Item(20,8).Value  = 0
$ws.Cells.Item(20,10).Value = 26
$ws.Cells.Item(20,11).Value = -65.384615384615
$ws.Cells.Item(20,12).Value = -25
$ws.Cells.Item(20,14).Value = -93.75

# ---------------------------------------------------------------------------
# Row 21 - TOTAL (bold styles 17/18)
# ---------------------------------------------------------------------------
$ws.Cells.Item(21,3).Value  = 11
$ws.Cells.Item(21,4).Value  = 18
$ws.Cells.Item(21,5).Value  = -38.888888888888
$ws.Cells.Item(21,6).Value  = 33
$ws.Cells.Item(21,7).Value  = 49
$ws.Cells.Item(21,8).Value  = -32.653061224489
$ws.Cells.Item(21,9).Value  = 137
$ws.Cells.Item(21,10).Value = 208
$ws.Cells.Item(21,11).Value = -34.134615384615
$ws.Cells.Item(21,12).Value = -4.861111111111
$ws.Cells.Item(21,13).Value = -2.142857142857
$ws.Cells.Item(21,14).Value = -77.577741407528

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Cells.Item(22,12).Value = -50
$ws.Cells.Item(22,13).Value = -87.5

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Cells.Item(23,4).Value  = 1
$ws.Cells.Item(23,5).Value  = -100
$ws.Cells.Item(23,6).Value  = 3
$ws.Cells.Item(23,8).Value  = 50
$ws.Cells.Item(23,10).Value = 13
$ws.Cells.Item(23,11).Value = -30.769230769230
$ws.Cells.Item(23,13).Value = -35.714285714285

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Cells.Item(24,3).Value  = 18
$ws.Cells.Item(24,4).Value  = 4
$ws.Cells.Item(24,5).Value  = 350
$ws.Cells.Item(24,6).Value  = 59
$ws.Cells.Item(24,7).Value  = 19
$ws.Cells.Item(24,8).Value  = 210.526315789474
$ws.Cells.Item(24,9).Value  = 199
$ws.Cells.Item(24,10).Value = 161
$ws.Cells.Item(24,11).Value = 23.602484472049
$ws.Cells.Item(24,12).Value = 9.340659340659
$ws.Cells.Item(24,13).Value = 97.029702970297

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Cells.Item(25,3).Value  = 11
$ws.Cells.Item(25,4).Value  = 2
$ws.Cells.Item(25,5).Value  = 450
$ws.Cells.Item(25,6).Value  = 34
$ws.Cells.Item(25,7).Value  = 8
$ws.Cells.Item(25,8).Value  = 325
$ws.Cells.Item(25,9).Value  = 107
$ws.Cells.Item(25,10).Value = 74
$ws.Cells.Item(25,11).Value = 44.594594594594
$ws.Cells.Item(25,12).Value = -4.464285714285

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Cells.Item(26,3).Value  = 3
$ws.Cells.Item(26,4).Value  = 4
$ws.Cells.Item(26,5).Value  = -25
$ws.Cells.Item(26,6).Value  = 18
$ws.Cells.Item(26,7).Value  = 20
$ws.Cells.Item(26,8).Value  = -10
$ws.Cells.Item(26,9).Value  = 77
$ws.Cells.Item(26,10).Value = 100
$ws.Cells.Item(26,11).Value = -23
$ws.Cells.Item(26,12).Value = -10.465116279069
$ws.Cells.Item(26,13).Value = -47.972972972973

# ---------------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Cells.Item(27,3).Value  = 1
$ws.Cells.Item(27,6).Value  = 1
$ws.Cells.Item(27,8).Value  = 0
$ws.Cells.Item(27,9).Value  = 2
$ws.Cells.Item(27,11).Value = -33.333333333333
$ws.Cells.Item(27,12).Value = 0

# ---------------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------------
$ws.Cells.Item(29,4).Value  = 3
$ws.Cells.Item(29,5).Value  = -100
$ws.Cells.Item(29,7).Value  = 4
$ws.Cells.Item(29,10).Value = 6

# ---------------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------------
$ws.Cells.Item(30,4).Value  = 1
$ws.Cells.Item(30,5).Value  = -100
$ws.Cells.Item(30,7).Value  = 2
$ws.Cells.Item(30,10).Value = 4
